# Adds a new weekly price record for "Coco" (Mercado Mayorista Lo Valledor
# de Santiago) by inserting a new row above the existing row 68, which
# shifts every subsequent record down by one row (old row 87 becomes the
# new row 88). The newly inserted row 68 is then filled with the new
# observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68 - this shifts rows 68:87 down to 69:88
# and carries the row-68 formatting (e.g. the date style on column D) down
# with it, exactly like Excel's own Rows.Insert.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly record.
$ws.Range("A68").Value = 6
$ws.Range("B68").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C68").Value = "Metropolitana"
$ws.Range("D68").Value = 44841
$ws.Range("E68").Value = 13
$ws.Range("F68").Value = "Fruta"
$ws.Range("G68").Value = 100108
$ws.Range("H68").Value = "Tropicales y subtropicales"
$ws.Range("I68").Value = 100108007
$ws.Range("J68").Value = "Coco"
$ws.Range("K68").Value = "Sin especificar"
$ws.Range("L68").Value = "Primera"
$ws.Range("M68").Value = 150
$ws.Range("N68").Value = 28000
$ws.Range("O68").Value = 30000
$ws.Range("P68").Value = 29000
$ws.Range("Q68").Value = "`$/malla 20 unidades"
$ws.Range("R68").Value = "Perú"
$ws.Range("S68").Value = 1450
$ws.Range("T68").Value = 20
